# Generate Report for Handback
# Updates row 7 ("7075403a-35ae-4efc-9228-e3c252fd9888.md") on both the
# zh-cn and de-de detail sheets: the handback for this file has now come in,
# but it is not the latest version, so:
#   - Latest Target File / Latest Handback File (I/J) get filled in
#   - Latest Handback DateTime (K) gets a real timestamp
#   - Error Detail (P) gets the "not latest" warning message
# and a hyperlink is added on the new I7 value, mirroring the one on A7.

$wb = $excel.ActiveWorkbook

$warningMessage = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2f71b642ee7edbce9f63cbe867d591cf9f604252/e2e/7075403a-35ae-4efc-9228-e3c252fd9888.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dc3764fd9e02d206cbe0ba4704495e6fbd7f2022/e2e/7075403a-35ae-4efc-9228-e3c252fd9888.md."
$targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dc3764fd9e02d206cbe0ba4704495e6fbd7f2022/e2e/7075403a-35ae-4efc-9228-e3c252fd9888.md"
$displayName = "7075403a-35ae-4efc-9228-e3c252fd9888.md"

# ---------- zh-cn sheet ----------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("J7").Value = "7075403a-35ae-4efc-9228-e3c252fd9888.3d63f2f3af2466527a8acd3f86899c48d4960bc3.zh-cn.xlf"
$ws.Range("K7").Value = "2016-08-27 02:54:10"
$ws.Range("P7").Value = $warningMessage

$ws.Range("I7").Value = $displayName
$ws.Hyperlinks.Add($ws.Range("I7"), $targetUrl, "", "", $displayName)
$ws.Range("I7").Font.Underline = $true
$ws.Range("I7").Font.Color = 15570276

# ---------- de-de sheet ----------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("J7").Value = "7075403a-35ae-4efc-9228-e3c252fd9888.3d63f2f3af2466527a8acd3f86899c48d4960bc3.de-de.xlf"
$ws.Range("K7").Value = "2016-08-27 02:54:16"
$ws.Range("P7").Value = $warningMessage

$ws.Range("I7").Value = $displayName
$ws.Hyperlinks.Add($ws.Range("I7"), $targetUrl, "", "", $displayName)
$ws.Range("I7").Font.Underline = $true
$ws.Range("I7").Font.Color = 15570276
